$wb = $excel.ActiveWorkbook

# Suppress the "Delete sheet" confirmation prompt
$excel.DisplayAlerts = $false

# Remove the "DATA" worksheet
$wsData = $wb.Worksheets.Item("DATA")
$wsData.Delete() | Out-Null

# Remove the "Adjusted Expenditure" worksheet
$wsAdj = $wb.Worksheets.Item("Adjusted Expenditure")
$wsAdj.Delete() | Out-Null

$excel.DisplayAlerts = $true
